# location-1: add a new livestream-camera row (Balneario Camboriu, Brazil)
# at the bottom of the table (row 38), following the same layout/format as
# the existing rows, and recompute the "Status" helper formula for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clone the formatting of the last existing data row (37) onto the
#        new row (38) before writing any values, so new cells pick up the
#        same borders/fills as the rest of the table (style index 6 for
#        A/C/D/E, style index 3 for the Status column G).
$ws.Range("A37:E37").Copy()
$ws.Range("A38:E38").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$ws.Range("G37").Copy()
$ws.Range("G38").PasteSpecial(-4122)       # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

# --- 2) Fill in the row's data, column by column, in the same order the
#        values were originally entered (YouTube link, Country,
#        lat/long, Location, City, Category) so brand-new text lands in
#        the shared-string table in that order.
$ws.Range("F38").Value = "vz_sIkkAv7Y"
$ws.Range("E38").Value = "Brazil"
# Lat/long starts with "-": lead with an apostrophe so it is stored as
# text (quote-prefixed), matching how the author entered it.
$ws.Range("B38").Value = "'-27.005150441765377, -48.63292964378517"
$ws.Range("C38").Value = "5ª AVENIDA - VILA REAL - BALNEÁRIO CAMBORIÚ - SC - BC AO VIVO"
$ws.Range("D38").Value = "Balneário Camboriú"
$ws.Range("A38").Value = "LIVE, TRAFFIC"

# --- 3) Status formula, same as every other row in the table.
$ws.Range("G38").Formula = "=IsYouTubeVideoValid(F38)"

# --- 4) Leave the sheet scrolled/selected the way the author left it
#        after finishing data entry.
[void]$excel.ActiveWindow.ScrollRow
$excel.ActiveWindow.ScrollRow = 16
[void]$ws.Range("F42").Select()
